$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.3558140695095062, 0.8999999761581421),
    @(0.2877398729324341, 0.9150000214576721),
    @(0.316961258649826,  0.8974999785423279),
    @(0.2832776308059692, 0.9200000166893005),
    @(0.3423527479171753, 0.8812500238418579)
)

$startRow = 9
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$wb.Save()
